$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,7).Value = 19.95578266666667
$ws.Cells.Item(2,8).Value = 59.867348
$ws.Cells.Item(2,9).Value = 0.0117373419656925
$ws.Cells.Item(2,10).Value = 0.0117373419656925
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,13).Value = 30.99161333333333
$ws.Cells.Item(2,14).Value = 92.97484
$ws.Cells.Item(2,15).Value = 0.3599121977633812
$ws.Cells.Item(2,16).Value = 0.3599121977633811
$ws.Cells.Item(2,17).Value = 618.461900169369
$ws.Cells.Item(2,18).Value = 5566.15710152432
$ws.Cells.Item(2,19).Value = 0.004224412542772752
$ws.Cells.Item(2,20).Value = 0.004224412542772751
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,7).Value = 19.95578266666667
$ws.Cells.Item(3,8).Value = 59.867348
$ws.Cells.Item(3,9).Value = 0.0117373419656925
$ws.Cells.Item(3,10).Value = 0.0117373419656925
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,13).Value = 29.913269
$ws.Cells.Item(3,14).Value = 89.739807
$ws.Cells.Item(3,15).Value = 0.3473891556493311
$ws.Cells.Item(3,16).Value = 0.3473891556493311
$ws.Cells.Item(3,17).Value = 596.9426950135373
$ws.Cells.Item(3,18).Value = 5372.484255121836
$ws.Cells.Item(3,19).Value = 0.004077425315029377
$ws.Cells.Item(3,20).Value = 0.004077425315029377
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,7).Value = 19.95578266666667
$ws.Cells.Item(4,8).Value = 59.867348
$ws.Cells.Item(4,9).Value = 0.0117373419656925
$ws.Cells.Item(4,10).Value = 0.0117373419656925
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,13).Value = 19.150218
$ws.Cells.Item(4,14).Value = 57.450654
$ws.Cells.Item(4,15).Value = 0.2223955550134164
$ws.Cells.Item(4,16).Value = 0.2223955550134163
$ws.Cells.Item(4,17).Value = 382.157588427288
$ws.Cells.Item(4,18).Value = 3439.418295845592
$ws.Cells.Item(4,19).Value = 0.002610332680842447
$ws.Cells.Item(4,20).Value = 0.002610332680842446
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,7).Value = 19.95578266666667
$ws.Cells.Item(5,8).Value = 59.867348
$ws.Cells.Item(5,9).Value = 0.0117373419656925
$ws.Cells.Item(5,10).Value = 0.0117373419656925
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,13).Value = 6.053716000000001
$ws.Cells.Item(5,14).Value = 18.161148
$ws.Cells.Item(5,15).Value = 0.07030309157387134
$ws.Cells.Item(5,16).Value = 0.07030309157387132
$ws.Cells.Item(5,17).Value = 120.8066408217227
$ws.Cells.Item(5,18).Value = 1087.259767395504
$ws.Cells.Item(5,19).Value = 0.0008251714270479226
$ws.Cells.Item(5,20).Value = 0.0008251714270479224
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,7).Value = 1637.343343333333
$ws.Cells.Item(6,8).Value = 4912.03003
$ws.Cells.Item(6,9).Value = 0.9630320723052701
$ws.Cells.Item(6,10).Value = 0.9630320723052702
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,13).Value = 30.99161333333333
$ws.Cells.Item(6,14).Value = 92.97484
$ws.Cells.Item(6,15).Value = 0.3599121977633812
$ws.Cells.Item(6,16).Value = 0.3599121977633811
$ws.Cells.Item(6,17).Value = 50743.9117904939
$ws.Cells.Item(6,18).Value = 456695.2061144452
$ws.Cells.Item(6,19).Value = 0.3466069896600132
$ws.Cells.Item(6,20).Value = 0.3466069896600132
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,7).Value = 1637.343343333333
$ws.Cells.Item(7,8).Value = 4912.03003
$ws.Cells.Item(7,9).Value = 0.9630320723052701
$ws.Cells.Item(7,10).Value = 0.9630320723052702
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,13).Value = 29.913269
$ws.Cells.Item(7,14).Value = 89.739807
$ws.Cells.Item(7,15).Value = 0.3473891556493311
$ws.Cells.Item(7,16).Value = 0.3473891556493311
$ws.Cells.Item(7,17).Value = 48978.29187448935
$ws.Cells.Item(7,18).Value = 440804.6268704042
$ws.Cells.Item(7,19).Value = 0.3345468984613534
$ws.Cells.Item(7,20).Value = 0.3345468984613534
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,7).Value = 1637.343343333333
$ws.Cells.Item(8,8).Value = 4912.03003
$ws.Cells.Item(8,9).Value = 0.9630320723052701
$ws.Cells.Item(8,10).Value = 0.9630320723052702
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,13).Value = 19.150218
$ws.Cells.Item(8,14).Value = 57.450654
$ws.Cells.Item(8,15).Value = 0.2223955550134164
$ws.Cells.Item(8,16).Value = 0.2223955550134163
$ws.Cells.Item(8,17).Value = 31355.48196568217
$ws.Cells.Item(8,18).Value = 282199.3376911396
$ws.Cells.Item(8,19).Value = 0.2141740522160511
$ws.Cells.Item(8,20).Value = 0.2141740522160511
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,7).Value = 1637.343343333333
$ws.Cells.Item(9,8).Value = 4912.03003
$ws.Cells.Item(9,9).Value = 0.9630320723052701
$ws.Cells.Item(9,10).Value = 0.9630320723052702
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,13).Value = 6.053716000000001
$ws.Cells.Item(9,14).Value = 18.161148
$ws.Cells.Item(9,15).Value = 0.07030309157387134
$ws.Cells.Item(9,16).Value = 0.07030309157387132
$ws.Cells.Item(9,17).Value = 9912.011595030493
$ws.Cells.Item(9,18).Value = 89208.10435527444
$ws.Cells.Item(9,19).Value = 0.06770413196785249
$ws.Cells.Item(9,20).Value = 0.06770413196785248
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,7).Value = 17.50081933333334
$ws.Cells.Item(10,8).Value = 52.502458
$ws.Cells.Item(10,9).Value = 0.01029341242216722
$ws.Cells.Item(10,10).Value = 0.01029341242216722
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,13).Value = 30.99161333333333
$ws.Cells.Item(10,14).Value = 92.97484
$ws.Cells.Item(10,15).Value = 0.3599121977633812
$ws.Cells.Item(10,16).Value = 0.3599121977633811
$ws.Cells.Item(10,17).Value = 542.3786257951912
$ws.Cells.Item(10,18).Value = 4881.407632156721
$ws.Cells.Item(10,19).Value = 0.003704724687347093
$ws.Cells.Item(10,20).Value = 0.003704724687347092
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,7).Value = 17.50081933333334
$ws.Cells.Item(11,8).Value = 52.502458
$ws.Cells.Item(11,9).Value = 0.01029341242216722
$ws.Cells.Item(11,10).Value = 0.01029341242216722
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,13).Value = 29.913269
$ws.Cells.Item(11,14).Value = 89.739807
$ws.Cells.Item(11,15).Value = 0.3473891556493311
$ws.Cells.Item(11,16).Value = 0.3473891556493311
$ws.Cells.Item(11,17).Value = 523.5067164384008
$ws.Cells.Item(11,18).Value = 4711.560447945606
$ws.Cells.Item(11,19).Value = 0.003575819850087006
$ws.Cells.Item(11,20).Value = 0.003575819850087006
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,7).Value = 17.50081933333334
$ws.Cells.Item(12,8).Value = 52.502458
$ws.Cells.Item(12,9).Value = 0.01029341242216722
$ws.Cells.Item(12,10).Value = 0.01029341242216722
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,13).Value = 19.150218
$ws.Cells.Item(12,14).Value = 57.450654
$ws.Cells.Item(12,15).Value = 0.2223955550134164
$ws.Cells.Item(12,16).Value = 0.2223955550134163
$ws.Cells.Item(12,17).Value = 335.144505411948
$ws.Cells.Item(12,18).Value = 3016.300548707532
$ws.Cells.Item(12,19).Value = 0.002289209168609873
$ws.Cells.Item(12,20).Value = 0.002289209168609873
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,7).Value = 17.50081933333334
$ws.Cells.Item(13,8).Value = 52.502458
$ws.Cells.Item(13,9).Value = 0.01029341242216722
$ws.Cells.Item(13,10).Value = 0.01029341242216722
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,13).Value = 6.053716000000001
$ws.Cells.Item(13,14).Value = 18.161148
$ws.Cells.Item(13,15).Value = 0.07030309157387134
$ws.Cells.Item(13,16).Value = 0.07030309157387132
$ws.Cells.Item(13,17).Value = 105.9449900113094
$ws.Cells.Item(13,18).Value = 953.5049101017842
$ws.Cells.Item(13,19).Value = 0.0007236587161232467
$ws.Cells.Item(13,20).Value = 0.0007236587161232467
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,7).Value = 25.39612333333333
$ws.Cells.Item(14,8).Value = 76.18836999999999
$ws.Cells.Item(14,9).Value = 0.01493717330687017
$ws.Cells.Item(14,10).Value = 0.01493717330687017
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,13).Value = 30.99161333333333
$ws.Cells.Item(14,14).Value = 92.97484
$ws.Cells.Item(14,15).Value = 0.3599121977633812
$ws.Cells.Item(14,16).Value = 0.3599121977633811
$ws.Cells.Item(14,17).Value = 787.0668345123111
$ws.Cells.Item(14,18).Value = 7083.601510610799
$ws.Cells.Item(14,19).Value = 0.005376070873248155
$ws.Cells.Item(14,20).Value = 0.005376070873248154
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,7).Value = 25.39612333333333
$ws.Cells.Item(15,8).Value = 76.18836999999999
$ws.Cells.Item(15,9).Value = 0.01493717330687017
$ws.Cells.Item(15,10).Value = 0.01493717330687017
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,13).Value = 29.913269
$ws.Cells.Item(15,14).Value = 89.739807
$ws.Cells.Item(15,15).Value = 0.3473891556493311
$ws.Cells.Item(15,16).Value = 0.3473891556493311
$ws.Cells.Item(15,17).Value = 759.6810688271767
$ws.Cells.Item(15,18).Value = 6837.129619444589
$ws.Cells.Item(15,19).Value = 0.005189012022861355
$ws.Cells.Item(15,20).Value = 0.005189012022861354
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,7).Value = 25.39612333333333
$ws.Cells.Item(16,8).Value = 76.18836999999999
$ws.Cells.Item(16,9).Value = 0.01493717330687017
$ws.Cells.Item(16,10).Value = 0.01493717330687017
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,13).Value = 19.150218
$ws.Cells.Item(16,14).Value = 57.450654
$ws.Cells.Item(16,15).Value = 0.2223955550134164
$ws.Cells.Item(16,16).Value = 0.2223955550134163
$ws.Cells.Item(16,17).Value = 486.3412981882199
$ws.Cells.Item(16,18).Value = 4377.071683693979
$ws.Cells.Item(16,19).Value = 0.00332196094791298
$ws.Cells.Item(16,20).Value = 0.003321960947912979
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,7).Value = 25.39612333333333
$ws.Cells.Item(17,8).Value = 76.18836999999999
$ws.Cells.Item(17,9).Value = 0.01493717330687017
$ws.Cells.Item(17,10).Value = 0.01493717330687017
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,13).Value = 6.053716000000001
$ws.Cells.Item(17,14).Value = 18.161148
$ws.Cells.Item(17,15).Value = 0.07030309157387134
$ws.Cells.Item(17,16).Value = 0.07030309157387132
$ws.Cells.Item(17,17).Value = 153.7409181609733
$ws.Cells.Item(17,18).Value = 1383.66826344876
$ws.Cells.Item(17,19).Value = 0.00105012946284768
$ws.Cells.Item(17,20).Value = 0.00105012946284768
